$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "8,10,12,14_*" block down by 64 rows (to rows 65-125)
$ws.Rows("1:64").Insert()

# Re-create the title-row merges on the (currently blank) rows 1-61 BEFORE
# writing any values/formats into them, so Excel does not backfill the
# other merged cells with placeholder style cells.
$ws.Range("B1:E1").Merge()
$ws.Range("B9:E9").Merge()
$ws.Range("B17:E17").Merge()
$ws.Range("B25:E25").Merge()
$ws.Range("B33:E33").Merge()
$ws.Range("B41:E41").Merge()
$ws.Range("B49:E49").Merge()
$ws.Range("B57:E57").Merge()

# --- Block starting at row 1 (format copied from old row 65) ---
$ws.Range("B65").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "0,2,4,6_randread_4k"

$ws.Range("B66:E66").Copy()
$ws.Range("B2:E2").PasteSpecial(-4122)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4

$ws.Range("A67:E67").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)
$ws.Range("A3").Value = "IOPS"
$ws.Range("B3").Value = 6268
$ws.Range("C3").Value = 13000
$ws.Range("D3").Value = 19200
$ws.Range("E3").Value = 25000
$ws.Range("A68:E68").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Range("A4").Value = "BW(MB/s)"
$ws.Range("B4").Value = 25.7
$ws.Range("C4").Value = 53.2
$ws.Range("D4").Value = 78.4
$ws.Range("E4").Value = 103
$ws.Range("A69:E69").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A5").Value = "lat_avg"
$ws.Range("B5").Value = 159.08873
$ws.Range("C5").Value = 151.19658
$ws.Range("D5").Value = 154.2982
$ws.Range("E5").Value = 155.13049

# --- Block starting at row 9 (format copied from old row 73) ---
$ws.Range("B73").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "0,2,4,6_randread_128k"

$ws.Range("B74:E74").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 4

$ws.Range("A75:E75").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A11").Value = "IOPS"
$ws.Range("B11").Value = 3271
$ws.Range("C11").Value = 5610
$ws.Range("D11").Value = 7511
$ws.Range("E11").Value = 10200
$ws.Range("A76:E76").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A12").Value = "BW(MB/s)"
$ws.Range("B12").Value = 429
$ws.Range("C12").Value = 735
$ws.Range("D12").Value = 984
$ws.Range("E12").Value = 1342
$ws.Range("A77:E77").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A13").Value = "lat_avg"
$ws.Range("B13").Value = 297.54
$ws.Range("C13").Value = 347.44
$ws.Range("D13").Value = 389.63
$ws.Range("E13").Value = 377.76

# --- Block starting at row 17 (format copied from old row 81) ---
$ws.Range("B81").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "0,2,4,6_randwrite_4k"

$ws.Range("B82:E82").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 4

$ws.Range("A83:E83").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A19").Value = "IOPS"
$ws.Range("B19").Value = 19700
$ws.Range("C19").Value = 11200
$ws.Range("D19").Value = 15400
$ws.Range("E19").Value = 61600
$ws.Range("A84:E84").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A20").Value = "BW(MB/s)"
$ws.Range("B20").Value = 80.7
$ws.Range("C20").Value = 45.9
$ws.Range("D20").Value = 63
$ws.Range("E20").Value = 252
$ws.Range("A85:E85").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A21").Value = "lat_avg"
$ws.Range("B21").Value = 44.49
$ws.Range("C21").Value = 102.23
$ws.Range("D21").Value = 80.39
$ws.Range("E21").Value = 32.43

# --- Block starting at row 25 (format copied from old row 89) ---
$ws.Range("B89").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = "0,2,4,6_randwrite_128k"

$ws.Range("B90:E90").Copy()
$ws.Range("B26:E26").PasteSpecial(-4122)
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 4

$ws.Range("A91:E91").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A27").Value = "IOPS"
$ws.Range("B27").Value = 1179
$ws.Range("C27").Value = 2398
$ws.Range("D27").Value = 2752
$ws.Range("E27").Value = 2864
$ws.Range("A92:E92").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A28").Value = "BW(MB/s)"
$ws.Range("B28").Value = 155
$ws.Range("C28").Value = 314
$ws.Range("D28").Value = 361
$ws.Range("E28").Value = 375
$ws.Range("A93:E93").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A29").Value = "lat_avg"
$ws.Range("B29").Value = 493.49
$ws.Range("C29").Value = 399.62
$ws.Range("D29").Value = 455.96
$ws.Range("E29").Value = 364.17

# --- Block starting at row 33 (format copied from old row 97) ---
$ws.Range("B97").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("B33").Value = "0,2,4,6_read_4k"

$ws.Range("B98:E98").Copy()
$ws.Range("B34:E34").PasteSpecial(-4122)
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = 4

$ws.Range("A99:E99").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)
$ws.Range("A35").Value = "IOPS"
$ws.Range("B35").Value = 116000
$ws.Range("C35").Value = 229000
$ws.Range("D35").Value = 349000
$ws.Range("E35").Value = 449000
$ws.Range("A100:E100").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Range("A36").Value = "BW(MB/s)"
$ws.Range("B36").Value = 474
$ws.Range("C36").Value = 939
$ws.Range("D36").Value = 1428
$ws.Range("E36").Value = 1839
$ws.Range("A101:E101").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)
$ws.Range("A37").Value = "lat_avg"
$ws.Range("B37").Value = 8.378
$ws.Range("C37").Value = 8.37676
$ws.Range("D37").Value = 8.25699
$ws.Range("E37").Value = 8.5496

# --- Block starting at row 41 (format copied from old row 105) ---
$ws.Range("B105").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("B41").Value = "0,2,4,6_read_128k"

$ws.Range("B106:E106").Copy()
$ws.Range("B42:E42").PasteSpecial(-4122)
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = 3
$ws.Range("E42").Value = 4

$ws.Range("A107:E107").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)
$ws.Range("A43").Value = "IOPS"
$ws.Range("B43").Value = 3580
$ws.Range("C43").Value = 7160
$ws.Range("D43").Value = 11900
$ws.Range("E43").Value = 14900
$ws.Range("A108:E108").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)
$ws.Range("A44").Value = "BW(MB/s)"
$ws.Range("B44").Value = 469
$ws.Range("C44").Value = 939
$ws.Range("D44").Value = 1561
$ws.Range("E44").Value = 1952
$ws.Range("A109:E109").Copy()
$ws.Range("A45:E45").PasteSpecial(-4122)
$ws.Range("A45").Value = "lat_avg"
$ws.Range("B45").Value = 272.26
$ws.Range("C45").Value = 273.29
$ws.Range("D45").Value = 245.7
$ws.Range("E45").Value = 256.77

# --- Block starting at row 49 (format copied from old row 113) ---
$ws.Range("B113").Copy()
$ws.Range("B49").PasteSpecial(-4122)
$ws.Range("B49").Value = "0,2,4,6_write_4k"

$ws.Range("B114:E114").Copy()
$ws.Range("B50:E50").PasteSpecial(-4122)
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = 2
$ws.Range("D50").Value = 3
$ws.Range("E50").Value = 4

$ws.Range("A115:E115").Copy()
$ws.Range("A51:E51").PasteSpecial(-4122)
$ws.Range("A51").Value = "IOPS"
$ws.Range("B51").Value = 4476
$ws.Range("C51").Value = 8937
$ws.Range("D51").Value = 13500
$ws.Range("E51").Value = 17800
$ws.Range("A116:E116").Copy()
$ws.Range("A52:E52").PasteSpecial(-4122)
$ws.Range("A52").Value = "BW(MB/s)"
$ws.Range("B52").Value = 18.3
$ws.Range("C52").Value = 36.6
$ws.Range("D52").Value = 55.3
$ws.Range("E52").Value = 73.1
$ws.Range("A117:E117").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)
$ws.Range("A53").Value = "lat_avg"
$ws.Range("B53").Value = 219.17
$ws.Range("C53").Value = 221.88
$ws.Range("D53").Value = 220.2
$ws.Range("E53").Value = 220.04

# --- Block starting at row 57 (format copied from old row 121) ---
$ws.Range("B121").Copy()
$ws.Range("B57").PasteSpecial(-4122)
$ws.Range("B57").Value = "0,2,4,6_write_128k"

$ws.Range("B122:E122").Copy()
$ws.Range("B58:E58").PasteSpecial(-4122)
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = 2
$ws.Range("D58").Value = 3
$ws.Range("E58").Value = 4

$ws.Range("A123:E123").Copy()
$ws.Range("A59:E59").PasteSpecial(-4122)
$ws.Range("A59").Value = "IOPS"
$ws.Range("B59").Value = 1412
$ws.Range("C59").Value = 2659
$ws.Range("D59").Value = 4272
$ws.Range("E59").Value = 5231
$ws.Range("A124:E124").Copy()
$ws.Range("A60:E60").PasteSpecial(-4122)
$ws.Range("A60").Value = "BW(MB/s)"
$ws.Range("B60").Value = 185
$ws.Range("C60").Value = 349
$ws.Range("D60").Value = 560
$ws.Range("E60").Value = 686
$ws.Range("A125:E125").Copy()
$ws.Range("A61:E61").PasteSpecial(-4122)
$ws.Range("A61").Value = "lat_avg"
$ws.Range("B61").Value = 531.39
$ws.Range("C61").Value = 537.82
$ws.Range("D61").Value = 535.64
$ws.Range("E61").Value = 533.97

$ws.Range("A1").Select()